$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.500631809234619
$ws.Range("B1").Value = 3.424493312835693
$ws.Range("C1").Value = 4.196808815002441
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 3.302035093307495
